$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 737, shifting rows 737:832 down to 738:833.
$ws.Rows(737).Insert()

# Populate the newly inserted row 737 with the new weekly record.
$ws.Range("A737").Value2 = 4
$ws.Range("B737").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C737").Value = "Los Lagos"
$ws.Range("D737").Value2 = 45127
$ws.Range("E737").Value2 = 10
$ws.Range("F737").Value = "Fruta"
$ws.Range("G737").Value2 = 100102
$ws.Range("H737").Value = "Cítricos"
$ws.Range("I737").Value2 = 100102005
$ws.Range("J737").Value = "Naranja"
$ws.Range("K737").Value = "Fukumoto"
$ws.Range("L737").Value = "Primera"
$ws.Range("M737").Value2 = 600
$ws.Range("N737").Value2 = 18000
$ws.Range("O737").Value2 = 18500
$ws.Range("P737").Value2 = 18250
$ws.Range("Q737").Value = "$/caja 15 kilos empedrada"
$ws.Range("R737").Value = "Región de O'Higgins"
$ws.Range("S737").Value2 = 1217
$ws.Range("T737").Value2 = 15
